# Updated cryptos list on Fri Jan 19 06:46:05 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.367.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "'2.467.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'311.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'94.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.56%  "
$ws.Range("D7").Value = "'0.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("D10").Value = "'33.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.22%  "
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "'7.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").Value = "'2.846.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").Value = "'2.477.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'14.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Value = "'41.332.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.37%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.91%  "
$ws.Range("D20").Value = "'0.0₃0923"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "'11.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.84%  "
$ws.Range("D22").Value = "'68.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").Value = "'237.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -5.97%  "
$ws.Range("D27").Value = "'24.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.08%  "
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("D29").Value = "'9.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.30%  "
$ws.Range("D30").Value = "'36.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").Value = "'151.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.08%  "
$ws.Range("D32").Value = "'5.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("E33").Value = "  -5.10%  "
$ws.Range("D34").Value = "'2.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "'0.0746"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").Value = "'16.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.10%  "
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").Value = "'4.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("E41").Value = "  -7.69%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'20.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.23%  "
$ws.Range("D44").Value = "'1.988.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'0.0286"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("E46").Value = "  -8.20%  "
$ws.Range("D47").Value = "'8.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.48%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'2.709.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'69.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'96.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.86%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'74.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.76%  "
